$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BDTPTUMCF")

# Rename existing "hydrogen" entry (row 24) to "hydrogen combustion turbine"
$ws.Range("A24").Value = "hydrogen combustion turbine"

# Add new row 25 for "hydrogen combined cycle"
$ws.Range("A25").Value = "hydrogen combined cycle"
$ws.Range("B25").Value = 1

# Apply the new formatting (black font color + vertically centered) to A24,
# matching the style introduced for the plant-type label cells. Build the
# style once on A24, then copy/paste the formatting onto A25 so both cells
# share a single new cell style instead of generating duplicate styles.
$r1 = $ws.Range("A24")
$r1.Font.Color = 0
$r1.VerticalAlignment = -4108

$r1.Copy() | Out-Null
$ws.Range("A25").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Make BDTPTUMCF the active sheet/tab, with D24 selected, as in the saved file
$ws.Activate() | Out-Null
$ws.Range("D24").Select() | Out-Null
